$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 2.3
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 12
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AI2").Value = 19
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 9
$ws.Range("AT2").Value = 2.63
